$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97; this pushes the existing rows 97..147
# down to 98..148, matching the target diff.
$ws.Rows(97).Insert()

# Populate the newly inserted row 97 with the new weekly record.
$ws.Cells.Item(97, 1).Value2  = 11
$ws.Cells.Item(97, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(97, 3).Value2  = "Bíobío"
$ws.Cells.Item(97, 4).Value2  = 45016
$ws.Cells.Item(97, 5).Value2  = 8
$ws.Cells.Item(97, 6).Value2  = 100112001
$ws.Cells.Item(97, 7).Value2  = "Berenjena"
$ws.Cells.Item(97, 8).Value2  = "Sin especificar"
$ws.Cells.Item(97, 9).Value2  = "Primera"
$ws.Cells.Item(97, 10).Value2 = 180
$ws.Cells.Item(97, 11).Value2 = 6500
$ws.Cells.Item(97, 12).Value2 = 7000
$ws.Cells.Item(97, 13).Value2 = 6722
$ws.Cells.Item(97, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(97, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(97, 16).Value2 = 112
$ws.Cells.Item(97, 17).Value2 = 60
$ws.Cells.Item(97, 18).Value2 = "Hortaliza"
